# Generate Report for Archive
#
# Two files ("1ead05b5-f298-4b58-a271-48cc54bc14e1.md" and
# "26af309a-e83d-4cde-81d8-e0727530b000.md") moved from "Ready for
# handoff" back to "In Translation" status. Update their Status cells on
# all three report sheets: the per-language tables ("zh-cn", "de-de")
# and the roll-up "Overview" sheet (which repeats each language's status
# in its own column).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B8").Value = "In Translation"
$wsOverview.Range("C8").Value = "In Translation"
$wsOverview.Range("B9").Value = "In Translation"
$wsOverview.Range("C9").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C8").Value = "In Translation"
$wsZhCn.Range("C9").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C8").Value = "In Translation"
$wsDeDe.Range("C9").Value = "In Translation"
